$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "Le premier élément qui sera possible de créer est un
# répertoire." -> "On pourra créer un répertoire."
# (Replace within run boundaries so we don't disturb other text.)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Le premier élément qu", $true, $false, $false, $false, $false, $true, 1, $false, "On pourra créer un répertoire", 2)
$d.Content.Find.Execute("i sera possible de créer est un répertoire. Le répertoire ", $true, $false, $false, $false, $false, $true, 1, $false, ". Le répertoire ", 2)

# ------------------------------------------------------------------
# Change 2a: insert new text right after "Sous un répertoire " and
# remove the _GoBack bookmark that used to sit at that spot.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$insertPos = $d.Range($goBack.Start, $goBack.Start)
$insertPos.InsertBefore("ou directement à la racine de son arborescence de bibliothèque, ")
$goBack.Delete()

# ------------------------------------------------------------------
# Change 2b: split the "options" paragraph so that a brand new
# paragraph describing the search feature is created, and re-anchor
# the _GoBack bookmark at the end of that new paragraph.
# ------------------------------------------------------------------
$optionsPara = $d.Paragraphs(8)
$optionsPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(9)
$newPara.Range.InsertBefore("L’utilisateur peut faire une recherche par titre de ses livres. Cela lui retourne les répertoires et/ou livres correspondants.#")

$marker = $d.Content
$marker.Find.Execute("#", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$marker.Bookmarks.Add("_GoBack")

$marker2 = $d.Content
$marker2.Find.Execute("#", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
